$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 2 (Tyrese Haliburton / Indiana Pacers) with row 14 (Stephen Curry / Golden State Warriors)
$ws.Range("A2").Value = "Stephen Curry"
$ws.Range("C2").Value = "Golden State Warriors"

$ws.Range("A14").Value = "Tyrese Haliburton"
$ws.Range("C14").Value = "Indiana Pacers"

# Swap row 10 (Jarrett Allen / Cleveland Cavaliers) with row 11 (Jalen Duren / Detroit Pistons)
$ws.Range("A10").Value = "Jalen Duren"
$ws.Range("C10").Value = "Detroit Pistons"

$ws.Range("A11").Value = "Jarrett Allen"
$ws.Range("C11").Value = "Cleveland Cavaliers"

# Rotate rows 16, 17, 18 up by one (16<-17, 17<-18, 18<-16), keep column B (Position) attached to the player
$ws.Range("A16").Value = "Kevin Durant"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Phoenix Suns"

$ws.Range("A17").Value = "OG Anunoby"
$ws.Range("B17").Value = "SF,PF"
$ws.Range("C17").Value = "New York Knicks"

$ws.Range("A18").Value = "Mark Williams"
$ws.Range("B18").Value = "C"
$ws.Range("C18").Value = "Charlotte Hornets"
